# Update pl_mw results for the 380 kV case (rows 2-25, columns D,E,F,G,H,O)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, then new values for columns D, E, F, G, H, O (column 4,5,6,7,8,15)
$updates = @(
    ,@(2, 0.003282768059786889, 0.4354855266416848, 0.4268455756567278, 0.3367752227115943, 0.3434920082864039, 1.288635611322178)
    ,@(3, 0.00311798581107503, 0.3796236440277312, 0.3933003553665202, 0.301200860300014, 0.3306581927846537, 1.186000399322808)
    ,@(4, 0.00301825549965784, 0.3454298728940586, 0.3729420886923549, 0.2794993337179221, 0.3229880749787242, 1.12370741712374)
    ,@(5, 0.002977983059397005, 0.3315199486717972, 0.3647054950681934, 0.2706907303260806, 0.3199150623406126, 1.098503372128306)
    ,@(6, 0.00297131824670771, 0.3292116194363075, 0.3633413997772692, 0.2692301621707429, 0.3194079646665671, 1.094329137056604)
    ,@(7, 0.003017710873157498, 0.3452421832879935, 0.3728307665163584, 0.2793803970866549, 0.3229464183999653, 1.123366776377708)
    ,@(8, 0.00322565360309568, 0.4162010468717909, 0.4152293778962388, 0.3244794999736769, 0.3390232442816625, 1.253095441307863)
    ,@(9, 0.003644731221282171, 0.5562976445451824, 0.500290933977297, 0.4140703990656505, 0.3722245405283218, 1.513334097580127)
    ,@(10, 0.003959326332179103, 0.6599630079396235, 0.5639956345127786, 0.4806440831151235, 0.3976543754028228, 1.708230485774493)
    ,@(11, 0.004103863068738889, 0.7073182102085269, 0.593248511635764, 0.5111043360736858, 0.4094518639923592, 1.797728296575144)
    ,@(12, 0.004158796916037488, 0.7252814645633947, 0.604365721277091, 0.5226648403902345, 0.4139524800982883, 1.831741466706035)
    ,@(13, 0.0041469570567223, 0.721411343475495, 0.6019696543234261, 0.5201739184789744, 0.4129817139520924, 1.824410655378472)
    ,@(14, 0.004108378505367938, 0.7087954249357438, 0.5941623311389179, 0.5120549037029605, 0.409821466670877, 1.800524119720649)
    ,@(15, 0.00408477407517438, 0.7010719049055893, 0.5893853151532795, 0.5070851600557091, 0.4078900483661414, 1.785908904299276)
    ,@(16, 0.003949908852849404, 0.6568724169760856, 0.5620894304974513, 0.4786570196495461, 0.3968880122874339, 1.702398612130537)
    ,@(17, 0.003867535620493356, 0.629809763319642, 0.5454146261147201, 0.4612626808594769, 0.3901974874650023, 1.651383722660228)
    ,@(18, 0.003820291157719424, 0.6142624613427188, 0.5358494138643692, 0.4512744032265061, 0.3863708487733675, 1.622120123892671)
    ,@(19, 0.003804318213450841, 0.6090015003014315, 0.5326151957723084, 0.4478953566152484, 0.3850789144428859, 1.612225456820568)
    ,@(20, 0.003876290507584201, 0.6326887045236305, 0.5471870247896504, 0.4631126271813173, 0.3909074710753089, 1.656806175146869)
    ,@(21, 0.004119704539757407, 0.7125001682737064, 0.5964544477322562, 0.5144389481118878, 0.4107488063637277, 1.807536841490787)
    ,@(22, 0.004279959021001645, 0.764842669832106, 0.6288856518110322, 0.5481347357245454, 0.423909655295347, 1.906761666784746)
    ,@(23, 0.00419432250913232, 0.7368890890761435, 0.6115551107154289, 0.5301366292088687, 0.4168677050772089, 1.85373763995716)
    ,@(24, 0.003872332070940843, 0.6313870999755409, 0.5463856568671019, 0.4622762292457878, 0.3905864257188227, 1.654354480355266)
    ,@(25, 0.003530170004236055, 0.5182800985497948, 0.477070003046336, 0.3897053291207442, 0.3630617902832967, 1.442293097308493)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 4).Value  = $u[1]   # D
    $ws.Cells.Item($r, 5).Value  = $u[2]   # E
    $ws.Cells.Item($r, 6).Value  = $u[3]   # F
    $ws.Cells.Item($r, 7).Value  = $u[4]   # G
    $ws.Cells.Item($r, 8).Value  = $u[5]   # H
    $ws.Cells.Item($r, 15).Value = $u[6]   # O
}

Write-Host "Updated $($updates.Count) rows (D,E,F,G,H,O) for the 380 kV case"